# SDKv5.00 update: replace the ONNX node-name labels in column B (rows 8-78)
# of the ResNet50 layer-summary sheet with their numeric IDs.
# Each new value is entered with a leading apostrophe so Excel stores it
# as text (matching the original inlineStr/text cell type) rather than
# auto-converting the numeric-looking string into a real number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B8"  = "496"
    "B9"  = "324"
    "B10" = "499"
    "B11" = "502"
    "B12" = "505"
    "B13" = "508"
    "B14" = "335"
    "B15" = "511"
    "B16" = "514"
    "B17" = "517"
    "B18" = "345"
    "B19" = "520"
    "B20" = "523"
    "B21" = "526"
    "B22" = "355"
    "B23" = "529"
    "B24" = "532"
    "B25" = "535"
    "B26" = "538"
    "B27" = "367"
    "B28" = "541"
    "B29" = "544"
    "B30" = "547"
    "B31" = "377"
    "B32" = "550"
    "B33" = "553"
    "B34" = "556"
    "B35" = "387"
    "B36" = "559"
    "B37" = "562"
    "B38" = "565"
    "B39" = "397"
    "B40" = "568"
    "B41" = "571"
    "B42" = "574"
    "B43" = "577"
    "B44" = "409"
    "B45" = "580"
    "B46" = "583"
    "B47" = "586"
    "B48" = "419"
    "B49" = "589"
    "B50" = "592"
    "B51" = "595"
    "B52" = "429"
    "B53" = "598"
    "B54" = "601"
    "B55" = "604"
    "B56" = "439"
    "B57" = "607"
    "B58" = "610"
    "B59" = "613"
    "B60" = "449"
    "B61" = "616"
    "B62" = "619"
    "B63" = "622"
    "B64" = "459"
    "B65" = "625"
    "B66" = "628"
    "B67" = "631"
    "B68" = "634"
    "B69" = "471"
    "B70" = "637"
    "B71" = "640"
    "B72" = "643"
    "B73" = "481"
    "B74" = "646"
    "B75" = "649"
    "B76" = "652"
    "B77" = "491"
    "B78" = "493"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = "'" + $updates[$addr]
}
